$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC32_Verify_store_location")

# Row 24 (VERIFY_WEBELEMENT_PRESENT / ContactUsConfirmation / CSS / ContactUsConfirmation)
# is an erroneous duplicate step that needs to be removed entirely, shifting
# the subsequent rows up by one.
$ws.Rows("24").Delete()

# Restore the selection/active cell as recorded after the edit.
$ws.Range("C21").Select()
